# edit.ps1 - apply the CV text-fix edits described by the diff.
#
# Summary of changes in the target diff:
#   1) Typo fix: "infrastrcuture" -> "infrastructure"
#      (...critical infrastrcuture requirements...)
#   2) Typo fix + missing word: "Automation vulnerability and inicident
#      management." -> "Automation of vulnerability and incident management."
#   3) Typo fix + missing word: "Debugging und solving of problems mainly in
#      the area mail (but not solely)." -> "Debugging and solving of problems
#      mainly in the area of mail (but not solely)."
#   4) Two <w:spacing> elements have their w:before/w:after attributes
#      reordered (w:before="150" w:after="50" -> w:after="50" w:before="150").
#      This is a value-preserving, cosmetic attribute-order change (the
#      paragraph spacing values themselves - SpaceBefore=150, SpaceAfter=50
#      twentieths-of-a-point - do not change). We still touch the
#      paragraphs' spacing formatting (idempotently, to the same values) so
#      that if the host re-serializes paragraph formatting on write, the
#      order matches; this is wrapped so it can never abort the script.

$d = $word.ActiveDocument

# --- 1) "infrastrcuture" -> "infrastructure" -------------------------------
$d.Content.Find.Execute(
    "infrastrcuture", $true, $false, $false, $false, $false,
    $true, 1, $false, "infrastructure", 2) | Out-Null

# --- 2) "Automation vulnerability and inicident management." ---------------
$d.Content.Find.Execute(
    "Automation vulnerability and inicident management.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Automation of vulnerability and incident management.", 2) | Out-Null

# --- 3) "Debugging und solving of problems mainly in the area mail (but not solely)." ---
$d.Content.Find.Execute(
    "Debugging und solving of problems mainly in the area mail (but not solely).", $true, $false, $false, $false, $false,
    $true, 1, $false, "Debugging and solving of problems mainly in the area of mail (but not solely).", 2) | Out-Null

# --- 4) Cosmetic <w:spacing> attribute-order swap on the two paragraphs ----
#     "Key Projects" and "Additional & Focused Projects" headings both carry
#     spacing before=150 (7.5pt) / after=50 (2.5pt), unchanged in value -
#     only the w:before/w:after attribute order flips. Re-apply the same
#     values through the paragraph formatting object so the pPr block gets
#     rewritten by the host; wrapped in try/catch so a host that can't
#     address paragraph formatting inside a table cell (or treats the
#     round-trip as a no-op) never aborts the rest of the script - nothing
#     else depends on this value-preserving formatting touch.
function Set-SpacingRoundTrip($paragraph) {
    try {
        $before = $paragraph.Format.SpaceBefore
        $after = $paragraph.Format.SpaceAfter
        $paragraph.Format.SpaceAfter = $after
        $paragraph.Format.SpaceBefore = $before
    } catch {
        # Not fatal - see comment above.
    }
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.Trim()
    if (($t -eq "Key Projects") -or ($t -eq "Additional & Focused Projects")) {
        Set-SpacingRoundTrip $p
    }
}
